# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.895.25"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "3.136.71"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'601.32"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").Value = "'143.09"
$ws.Range("E6").Value = "  -3.56%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.130.11"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "3.649.60"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "63.933.85"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "3.124.69"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").Value = "'485.05"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'14.66"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'0.708"
$ws.Range("E22").Value = "  -1.55%  "
$ws.Range("D23").Value = "'7.63"
$ws.Range("E23").Value = "  -4.77%  "
$ws.Range("D24").Value = "'86.91"
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("D25").Value = "'13.44"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("D28").Value = "'8.27"
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").Value = "'27.12"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.111"
$ws.Range("E32").Value = "  -7.85%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "'2.65"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").Value = "'6.00"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'52.52"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("E38").Value = "  -5.80%  "
$ws.Range("E39").Value = "  -7.81%  "
$ws.Range("D40").Value = "'438.45"
$ws.Range("D41").Value = "'0.0395"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "'8.28"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").Value = "2.876.61"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("D47").Value = "'2.38"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "'25.92"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'121.35"
$ws.Range("E51").Value = "  +0.86%  "
